# Lamella/SLD_SasView.xlsx update
# - removes Sheet2 (now redundant, data folded into/kept only on Sheet1)
# - refreshes SLD / dSLD header labels (Å**-2 -> Å^-2)
# - reworks the "q/d-spacing" helper block (rows 14-23) to reference a
#   live input cell (C16) instead of a hard-coded 0.1, and relabels the
#   SasView fit block

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Drop Sheet2 - its content was a duplicate/stale copy of Sheet1
# ---------------------------------------------------------------------
$excel.DisplayAlerts = $false
$ws2 = $wb.Worksheets.Item("Sheet2")
if ($ws2) {
    $ws2.Delete() | Out-Null
}

$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 2) Header row: "Å**-2" -> "Å^-2"
# ---------------------------------------------------------------------
$ws.Range("G1").Value = "SLD, E-6 Å^-2"
$ws.Range("H1").Value = "dSLD, E-6 Å^-2"

# ---------------------------------------------------------------------
# 3) Rework the q / d-spacing block
# ---------------------------------------------------------------------
# old row 17 (A:H) held: d-spacing | q= | 2pi/d | so  | d= | 2pi/q | = | =2*PI()/0.1
# clear it out first, then rebuild rows 14, 16, 17, 20, 23 as needed
$ws.Range("A17:H17").ClearContents()

$ws.Range("A14").Value = "SASVIEW MODELLING:"
$ws.Range("A14").Font.Bold = $true

$ws.Range("B16").Value = "q_first_peak ="
$ws.Range("C16").Value = 0.1

$ws.Range("A17").Value = "d-spacing"
$ws.Range("B17").Value = "d ="
$ws.Range("C17").Value = "2pi/q_first_peak ="
$ws.Range("D17").Formula = "=2*PI()/C16"
$ws.Range("D17").NumberFormat = "0"
$ws.Range("E17").Value = "Å"

$ws.Range("B20").Value = "fit sample 2"

$ws.Range("B23").Value = "SASVIEW fit sample 2"
$ws.Range("C23").Value = "SASVIEW fit sample 1"

# rows 24-29 (length_tail/length_head/sld_head/sld_tail/chi2r) are untouched

# ---------------------------------------------------------------------
# 4) Cosmetic: column B widened to fit the new "SASVIEW fit sample 2" text
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 18.45

# ---------------------------------------------------------------------
# 5) Restore the last-used selection
# ---------------------------------------------------------------------
$ws.Range("C24").Select() | Out-Null

$wb.Save() | Out-Null
